# Apply updated secure credential loader values to "Top 100" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M4").Value = 0.459
$ws.Range("R4").Value = 0.09299999999999999
$ws.Range("L5").Value = 41
$ws.Range("R5").Value = 0.1760869565217392
$ws.Range("N6").Value = 0.4033333333333334
$ws.Range("Q6").Value = 0.6791666666666667
$ws.Range("N7").Value = 0.501470588235294
$ws.Range("L9").Value = 30
$ws.Range("R9").Value = 0.3307142857142857
$ws.Range("P10").Value = 0.4392307692307693
$ws.Range("S10").Value = 0.2553846153846153
$ws.Range("Q11").Value = 0.6633333333333333
$ws.Range("Q12").Value = 0.67875
$ws.Range("Q13").Value = 0.8409090909090909
$ws.Range("N14").Value = 0.41
$ws.Range("O14").Value = 0.03333333333333335
$ws.Range("Q15").Value = 0.5508333333333333
$ws.Range("M17").Value = 0.6009999999999999
$ws.Range("O17").Value = 0.05950000000000002
$ws.Range("L19").Value = 73
$ws.Range("N19").Value = 0.5363636363636364
$ws.Range("Q19").Value = 0.4136363636363635
$ws.Range("R21").Value = 0.1320000000000001
$ws.Range("N22").Value = 0.6869230769230769
$ws.Range("P22").Value = 0.4876923076923078
$ws.Range("N23").Value = 0.6145454545454546
$ws.Range("P23").Value = 0.5027272727272727
$ws.Range("N24").Value = 0.403
$ws.Range("M25").Value = 0.3545454545454545
$ws.Range("N25").Value = 0.4540909090909091
$ws.Range("O25").Value = 0.0531818181818182
$ws.Range("P25").Value = 0.4390909090909092
$ws.Range("R25").Value = 0.8454545454545453
$ws.Range("N26").Value = 0.3317857142857142
$ws.Range("Q26").Value = 0.7957142857142859
$ws.Range("M27").Value = 0.5
$ws.Range("N28").Value = 0.5946666666666667
$ws.Range("P28").Value = 0.5473333333333332
$ws.Range("L29").Value = 45
$ws.Range("M29").Value = 0.64
$ws.Range("P29").Value = 0.2735714285714285
$ws.Range("L30").Value = 104
$ws.Range("P30").Value = 0.3391666666666666
$ws.Range("R30").Value = 0.2291666666666666
$ws.Range("S30").Value = 0.2999999999999999
$ws.Range("N32").Value = 0.5433333333333334
$ws.Range("O32").Value = 0.07433333333333336
$ws.Range("P32").Value = 0.2736666666666667
$ws.Range("Q32").Value = 0.5813333333333335
$ws.Range("R32").Value = 0.2473333333333333
$ws.Range("N35").Value = 0.5345454545454547
$ws.Range("N36").Value = 0.3754545454545455
$ws.Range("Q37").Value = 0.9662500000000001
$ws.Range("L38").Value = 7
$ws.Range("M38").Value = 0.4218750000000001
$ws.Range("Q39").Value = 0.8236363636363638
$ws.Range("L40").Value = 19
$ws.Range("P41").Value = 0.6927272727272729
$ws.Range("Q41").Value = 0.7836363636363637
$ws.Range("O42").Value = 0.04250000000000001
$ws.Range("P42").Value = 0.4125000000000001
$ws.Range("S42").Value = 0.4483333333333334
$ws.Range("Q44").Value = 0.7119999999999999
$ws.Range("Q46").Value = 0.5952380952380952
$ws.Range("L47").Value = 25
$ws.Range("M47").Value = 0.2033333333333333
$ws.Range("N47").Value = 0.268
$ws.Range("O47").Value = 0.05533333333333334
$ws.Range("P47").Value = 0.1366666666666667
$ws.Range("Q47").Value = 0.5533333333333332
$ws.Range("R47").Value = 0.16
$ws.Range("S47").Value = 0.5433333333333332
$ws.Range("R48").Value = 0.1358333333333334
$ws.Range("L51").Value = 22
$ws.Range("N51").Value = 0.3458333333333333
$ws.Range("Q51").Value = 0.8775000000000001
$ws.Range("R51").Value = 0.1566666666666666
$ws.Range("O52").Value = 0.0675
$ws.Range("S52").Value = 0.5075
$ws.Range("N53").Value = 0.6858823529411766
$ws.Range("N54").Value = 0.383888888888889
$ws.Range("R54").Value = 0.1844444444444443
$ws.Range("M55").Value = 0.0742857142857143
$ws.Range("N55").Value = 0.5735714285714286
$ws.Range("O55").Value = 0.05928571428571429
$ws.Range("N56").Value = 0.5099999999999999
$ws.Range("P56").Value = 0.4407692307692308
$ws.Range("M59").Value = 0.431875
$ws.Range("Q59").Value = 0.6031249999999999
$ws.Range("S60").Value = 0.61
$ws.Range("M61").Value = 0.03333333333333334
$ws.Range("N61").Value = 0.4233333333333334
$ws.Range("O61").Value = 0.04777777777777778
$ws.Range("P61").Value = 0.1533333333333333
$ws.Range("Q61").Value = 0.7255555555555555
$ws.Range("R61").Value = 0.1644444444444445
$ws.Range("S61").Value = 0.27
$ws.Range("M62").Value = 0.4383333333333332
$ws.Range("P62").Value = 0.1466666666666666
$ws.Range("Q62").Value = 0.5191666666666668
$ws.Range("R62").Value = 0.1333333333333333
$ws.Range("S63").Value = 0.6929999999999999
$ws.Range("Q64").Value = 0.9225
$ws.Range("M65").Value = 0.052
$ws.Range("N66").Value = 0.4961538461538461
$ws.Range("O66").Value = 0.09769230769230772
$ws.Range("Q66").Value = 0.8761538461538462
$ws.Range("R66").Value = 0.2507692307692308
$ws.Range("Q67").Value = 0.8303846153846153
$ws.Range("R67").Value = 0.1961538461538461
$ws.Range("P68").Value = 0.4819230769230769
$ws.Range("Q68").Value = 0.5765384615384617
$ws.Range("R68").Value = 0.2284615384615385
$ws.Range("Q69").Value = 0.4584615384615385
$ws.Range("R69").Value = 0.3176923076923078
$ws.Range("N71").Value = 0.535
$ws.Range("Q71").Value = 0.5200000000000001
$ws.Range("M72").Value = 0.662142857142857
$ws.Range("O72").Value = 0.04357142857142858
$ws.Range("N73").Value = 0.4687500000000001
$ws.Range("O73").Value = 0.04875
$ws.Range("Q73").Value = 0.7206249999999998
$ws.Range("Q74").Value = 0.9466666666666668
$ws.Range("N77").Value = 0.2441666666666667
$ws.Range("S77").Value = 0.3116666666666666
